$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value = 45
Write-Host "test ok"
